# Apply the edits described by the diff across the three worksheets:
# Summary, Assets, Liabilities.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: Summary
# -----------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Sama Al Qubaisi"
$wsSummary.Range("B4").Value = 5957
$wsSummary.Range("B6").Value = 258359
$wsSummary.Range("B7").Value = 201980
$wsSummary.Range("B8").Value = 56379
$wsSummary.Range("B9").Value = 1.28

# -----------------------------------------------------------------
# Sheet 2: Assets
#   Row 2 (Vehicles / Luxury Car)       -> Vehicles / Premium Car, 254019
#   Row 3 (Vehicles / Mid-range Car)    -> Liquid Assets / Savings Account, 4340
#   Row 4 (Liquid Assets / Savings Acc) -> removed (row shift up)
#   Row 5 (TOTAL ASSETS)                -> becomes row 4, 258359
# -----------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("B2").Value = "Premium Car"
$wsAssets.Range("C2").Value = 254019
$wsAssets.Range("A3").Value = "Liquid Assets"
$wsAssets.Range("B3").Value = "Savings Account"
$wsAssets.Range("C3").Value = 4340
$wsAssets.Rows(4).Delete()
$wsAssets.Range("C4").Value = 258359

# -----------------------------------------------------------------
# Sheet 3: Liabilities
#   Row 2 (Auto Loans / Vehicle Loan 1) -> 152411 / 2540 / 5
#   Row 3 (Auto Loans / Vehicle Loan 2) -> Personal Loans / Personal Loan, 33491 / 698 / 4
#   Row 4 (Personal Loans / Personal Loan) -> Credit Cards / Credit Card Balance, 16078 / 804 / 1
#   Row 5 (Credit Cards / Credit Card Balance) -> removed (row shift up)
#   Row 6 (TOTAL LIABILITIES)           -> becomes row 5, 201980
# -----------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Range("C2").Value = 152411
$wsLiabilities.Range("D2").Value = 2540
$wsLiabilities.Range("E2").Value = 5

$wsLiabilities.Range("A3").Value = "Personal Loans"
$wsLiabilities.Range("B3").Value = "Personal Loan"
$wsLiabilities.Range("C3").Value = 33491
$wsLiabilities.Range("D3").Value = 698
$wsLiabilities.Range("E3").Value = 4

$wsLiabilities.Range("A4").Value = "Credit Cards"
$wsLiabilities.Range("B4").Value = "Credit Card Balance"
$wsLiabilities.Range("C4").Value = 16078
$wsLiabilities.Range("D4").Value = 804
$wsLiabilities.Range("E4").Value = 1

$wsLiabilities.Rows(5).Delete()
$wsLiabilities.Range("C5").Value = 201980

Write-Host "Edits applied"
